$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.692.57"
$ws.Range("E2").Value = "  -0.72%  "
$ws.Range("D3").Value = "1.679.16"
$ws.Range("E3").Value = "  -1.27%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9988"
$ws.Range("E4").Value = "  -0.47%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.84"
$ws.Range("E5").Value = "  -0.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9997"
$ws.Range("E6").Value = "  -0.32%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3934"
$ws.Range("E7").Value = "  -1.85%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3959"
$ws.Range("E8").Value = "  -2.61%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.001"
$ws.Range("E9").Value = "  -0.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.410"
$ws.Range("E10").Value = "  -3.65%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "50.99"
$ws.Range("E11").Value = "  -5.17%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08676"
$ws.Range("E12").Value = "  -1.58%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "25.42"
$ws.Range("E13").Value = "  -1.51%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.356"
$ws.Range("E14").Value = "  -1.63%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001323"
$ws.Range("E15").Value = "  -1.56%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.737"
$ws.Range("E16").Value = "  -3.82%  "
$ws.Range("D17").Value = "1.671.34"
$ws.Range("E17").Value = "  -7.99%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.14"
$ws.Range("E18").Value = "  -2.60%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07037"
$ws.Range("E19").Value = "  -2.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "21.36"
$ws.Range("E20").Value = "  +1.85%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.100"
$ws.Range("E21").Value = "  -1.92%  "
$ws.Range("E22").Value = "  -0.30%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.00"
$ws.Range("E23").Value = "  -3.76%  "
$ws.Range("D24").Value = "24.652.55"
$ws.Range("E24").Value = "  -0.90%  "
$ws.Range("E25").Value = "  +1.26%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.791"
$ws.Range("E26").Value = "  -3.37%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.22"
$ws.Range("E27").Value = "  +0.57%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.870"
$ws.Range("E28").Value = "  -12.38%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "160.38"
$ws.Range("E29").Value = "  -1.80%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "146.89"
$ws.Range("E30").Value = "  +2.38%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.316"
$ws.Range("E31").Value = "  +1.72%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.521"
$ws.Range("E32").Value = "  +10.97%  "
$ws.Range("D33").Value = "1.853.01"
$ws.Range("E33").Value = "  -1.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.03097"
$ws.Range("E34").Value = "  -2.51%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08330"
$ws.Range("E35").Value = "  -4.52%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.987"
$ws.Range("E36").Value = "  -5.67%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2817"
$ws.Range("E37").Value = "  -2.23%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.9915"
$ws.Range("E38").Value = "  -4.28%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.09546"
$ws.Range("E39").Value = "  +1.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.515"
$ws.Range("E40").Value = "  +2.85%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "10.36"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7931"
$ws.Range("E42").Value = "  -7.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "13.62"
$ws.Range("E43").Value = "  -3.09%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.51"
$ws.Range("E44").Value = "  -6.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.7154"
$ws.Range("E45").Value = "  -4.23%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.571"
$ws.Range("E46").Value = "  -4.74%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.171"
$ws.Range("E47").Value = "  -1.33%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.08663"
$ws.Range("E48").Value = "  +3.77%  "
$ws.Range("E49").Value = "  -0.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.337"
$ws.Range("E50").Value = "  -4.94%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "137.80"
$ws.Range("E51").Value = "  -2.39%  "
